$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Cells.Item(3, 3).Value = 0.1999146332477239
$ws.Cells.Item(3, 4).Value = 0.9419246634644478
$ws.Cells.Item(3, 5).Value = 0.2763937264113785
$ws.Cells.Item(3, 7).Value = 'max\_depth: 2, max\_features: 6 \\'

# Row 4
$ws.Cells.Item(4, 3).Value = 0.2098865126194012
$ws.Cells.Item(4, 4).Value = 0.9889085133642006
$ws.Cells.Item(4, 5).Value = 0.2833951531294875

# Row 5
$ws.Cells.Item(5, 3).Value = 0.2235364021146156
$ws.Cells.Item(5, 4).Value = 1.053221802292754
$ws.Cells.Item(5, 5).Value = 0.2944560654084189

# Row 6
$ws.Cells.Item(6, 3).Value = 0.2078335327122884
$ws.Cells.Item(6, 4).Value = 0.9792356226073227
$ws.Cells.Item(6, 5).Value = 0.2833126417239996

# Row 7
$ws.Cells.Item(7, 3).Value = 0.2106557187265478
$ws.Cells.Item(7, 4).Value = 0.9925327313208274
$ws.Cells.Item(7, 5).Value = 0.2837120684656571
$ws.Cells.Item(7, 7).Value = 'learning\_rate: 0.05, max\_depth: 2, max\_features: 8, n\_estimators: 50 \\'

# Row 10
$ws.Cells.Item(10, 3).Value = 0.2115509307079892
$ws.Cells.Item(10, 4).Value = 0.9967506428896301
$ws.Cells.Item(10, 5).Value = 0.2961002766162041

# Row 11
$ws.Cells.Item(11, 3).Value = 0.2193257768158176
$ws.Cells.Item(11, 4).Value = 1.03338287528119
$ws.Cells.Item(11, 5).Value = 0.2952673527117101

# Row 12
$ws.Cells.Item(12, 3).Value = 0.2992442535360612
$ws.Cells.Item(12, 4).Value = 1.409929519548232
$ws.Cells.Item(12, 5).Value = 0.3728801803282877

# Row 13
$ws.Cells.Item(13, 3).Value = 0.2207354344330774
$ws.Cells.Item(13, 4).Value = 1.040024666605651
$ws.Cells.Item(13, 5).Value = 0.296385546190199

# Row 14
$ws.Cells.Item(14, 3).Value = 0.2617386551353075
$ws.Cells.Item(14, 4).Value = 1.233216851857284
$ws.Cells.Item(14, 5).Value = 0.3465569932215036

# Row 15
$ws.Cells.Item(15, 3).Value = 0.2169132860162558
$ws.Cells.Item(15, 4).Value = 1.022016100635572
$ws.Cells.Item(15, 5).Value = 0.296112048433617

# Row 16
$ws.Cells.Item(16, 3).Value = 0.2072415936425252
$ws.Cells.Item(16, 4).Value = 0.9764466221223643
$ws.Cells.Item(16, 5).Value = 0.2867524410915621
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 'omega: 50, r: 0.1 \\'

# Row 17
$ws.Cells.Item(17, 3).Value = 0.2072415936425252
$ws.Cells.Item(17, 4).Value = 0.9764466221223643
$ws.Cells.Item(17, 5).Value = 0.2867524410915621
$ws.Cells.Item(17, 6).Value = 1

# Row 18
$ws.Cells.Item(18, 3).Value = 0.210084587743563
$ws.Cells.Item(18, 4).Value = 0.9898417709333728
$ws.Cells.Item(18, 5).Value = 0.2902756803086696
$ws.Cells.Item(18, 6).Value = 2

# Row 19
$ws.Cells.Item(19, 3).Value = 0.2102132155796646
$ws.Cells.Item(19, 4).Value = 0.9904478182709984
$ws.Cells.Item(19, 5).Value = 0.2911033245649785
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 'alpha: 0.9, beta: 0.1, lambda1: 0.001, omega: 100, sigma: 0.1 \\'

# Row 20
$ws.Cells.Item(20, 3).Value = 0.3036376979424972
$ws.Cells.Item(20, 4).Value = 1.430629823356678
$ws.Cells.Item(20, 5).Value = 0.4015482312448757
$ws.Cells.Item(20, 6).Value = 374
$ws.Cells.Item(20, 7).Value = 'alpha: 0.01, lambda1: 0.5, omega: 10000, sigma: 0.001, w: 50 \\'

# Row 21
$ws.Cells.Item(21, 3).Value = 0.2088349711560764
$ws.Cells.Item(21, 4).Value = 0.983954034430516
$ws.Cells.Item(21, 5).Value = 0.2897794988032813
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 'alpha: 0.001, beta: 0.01, e\_utility: 0.05, lambda1: 0.75, omega: 100, pi: 0.5, sigma: 0.5 \\'

# Row 22
$ws.Cells.Item(22, 3).Value = 0.2256003631239277
$ws.Cells.Item(22, 4).Value = 1.062946431988526
$ws.Cells.Item(22, 5).Value = 0.2983880575498681
$ws.Cells.Item(22, 6).Value = 19
$ws.Cells.Item(22, 7).Value = 'alpha: 0.1, beta: 0.1, e\_utility: 0.05, lambda1: 0.001, sigma: 50 \\'

# Row 23
$ws.Cells.Item(23, 3).Value = 0.3560309811869754
$ws.Cells.Item(23, 4).Value = 1.6774878191228
$ws.Cells.Item(23, 5).Value = 0.4136158655623717
$ws.Cells.Item(23, 7).Value = 'fuzzy\_operator: min, rules: 2 \\'

# Row 24
$ws.Cells.Item(24, 3).Value = 0.2189218688832968
$ws.Cells.Item(24, 4).Value = 1.031479808771102
$ws.Cells.Item(24, 5).Value = 0.3083392160051394

# Row 25
$ws.Cells.Item(25, 3).Value = 0.2149236833916936
$ws.Cells.Item(25, 4).Value = 1.012641820463462
$ws.Cells.Item(25, 5).Value = 0.2801843194044386
$ws.Cells.Item(25, 6).Value = 14
$ws.Cells.Item(25, 7).Value = 'adaptive\_filter: wRLS, fuzzy\_operator: prod, rules: 14 \\'

# Row 26
$ws.Cells.Item(26, 3).Value = 0.2333700016774168
$ws.Cells.Item(26, 4).Value = 1.09955412828791
$ws.Cells.Item(26, 5).Value = 0.3156199595167128
$ws.Cells.Item(26, 6).Value = 5
$ws.Cells.Item(26, 7).Value = 'error\_metric: RMSE, fuzzy\_operator: prod, num\_generations: 10, num\_parents\_mating: 5, parallel\_processing: 10, rules: 5, sol\_per\_pop: 10 \\'

# Row 27
$ws.Cells.Item(27, 3).Value = 0.2179841230615975
$ws.Cells.Item(27, 4).Value = 1.027061493297293
$ws.Cells.Item(27, 5).Value = 0.3014258757197226
$ws.Cells.Item(27, 7).Value = 'adaptive\_filter: RLS, error\_metric: MAE, fuzzy\_operator: minmax, lambda1: 0.99, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 1, sol\_per\_pop: 5 \\'

# Row 28
$ws.Cells.Item(28, 3).Value = 0.2071346074600338
$ws.Cells.Item(28, 4).Value = 0.9759425423444039
$ws.Cells.Item(28, 5).Value = 0.284029599924307
$ws.Cells.Item(28, 6).Value = 13
$ws.Cells.Item(28, 7).Value = 'adaptive\_filter: wRLS, error\_metric: RMSE, fuzzy\_operator: max, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 13, sol\_per\_pop: 5 \\'

# Row 29
$ws.Cells.Item(29, 3).Value = 0.2724769316005904
$ws.Cells.Item(29, 4).Value = 1.283811684668829
$ws.Cells.Item(29, 5).Value = 0.3518008894852795

# Row 30
$ws.Cells.Item(30, 3).Value = 0.2055663958294065
$ws.Cells.Item(30, 4).Value = 0.9685537024759927
$ws.Cells.Item(30, 5).Value = 0.2802669775120018
$ws.Cells.Item(30, 7).Value = 'combination: mean, n\_estimators: 50 \\'

# Row 31
$ws.Cells.Item(31, 3).Value = 0.2062990781373664
$ws.Cells.Item(31, 4).Value = 0.9720058336438813
$ws.Cells.Item(31, 5).Value = 0.2801898032974084
